$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "dSF" (column F) values from repulled data
$ws.Range("F3").Value = -1
$ws.Range("F4").Value = 0
$ws.Range("F6").Value = -1
$ws.Range("F8").Value = -1
$ws.Range("F9").Value = -1
$ws.Range("F10").Value = -2
$ws.Range("F11").Value = 0
$ws.Range("F13").Value = 1
$ws.Range("F14").Value = -2
$ws.Range("F16").Value = -3
$ws.Range("F17").Value = -4
$ws.Range("F18").Value = 1
$ws.Range("F19").Value = -1
$ws.Range("F21").Value = 4
$ws.Range("F22").Value = -6
$ws.Range("F23").Value = 2
$ws.Range("F24").Value = -1
$ws.Range("F25").Value = 1
$ws.Range("F28").Value = 0
$ws.Range("F29").Value = -10
$ws.Range("F30").Value = 1
